$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row ---
# J1: "operator" -> "operator_name"
$ws.Range("J1").Value = "operator_name"

# K1: new header "operator_code", matching the header style used by A1:J1
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "operator_code"
$excel.CutCopyMode = $false

# --- Row 2 ---
$ws.Range("B2").Value = 45971
$ws.Range("C2").Value = "Shift 2"
$ws.Range("D2").Value = "7 AM"
$ws.Range("E2").Value = "6 AM"
$ws.Range("F2").Value = "dh"
$ws.Range("G2").Value = "Assembly"
$ws.Range("H2").Value = "Pre-Assembly"
$ws.Range("I2").Value = "Pre Line 1"
$ws.Range("J2").Value = "Worker A1"
# K2: numeric-looking code must stay text, like the source workbook
$ws.Range("K2").Value = "'401"
$ws.Range("K2").Style = "Normal"

# --- Row 3 ---
$ws.Range("B3").Value = 45971
$ws.Range("C3").Value = "Shift 1"
$ws.Range("D3").Value = "9 AM"
$ws.Range("E3").Value = "8 AM"
$ws.Range("F3").Value = "dh"
$ws.Range("G3").Value = "Assembly"
$ws.Range("H3").Value = "Pre-Assembly"
$ws.Range("I3").Value = "Pre Line 1"
$ws.Range("J3").Value = "Worker A1"
$ws.Range("K3").Value = "'401"
$ws.Range("K3").Style = "Normal"
